# Auto-generated edit script: updates cryptos list (prices & 1h volume %)
# and fixes a few swapped coin rows, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are plain decimal numbers (e.g. '557.08').
# Force those specific cells to Text format first so Excel doesn't
# auto-convert them from strings into numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "63.273.51"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "3.022.03"
$ws.Range("E3").Value = "  -2.93%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "557.08"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "156.31"
$ws.Range("E6").Value = "  -3.76%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.559"
$ws.Range("E8").Value = "  -4.62%  "
$ws.Range("D9").Value = "3.030.42"
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("D11").Value = "6.41"
$ws.Range("E11").Value = "  -4.52%  "
$ws.Range("D12").Value = "0.367"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "3.549.37"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "63.234.21"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "24.10"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "3.023.15"
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").Value = "0.0000151"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "396.91"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "5.11"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "12.06"
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "65.16"
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("E25").Value = "  -5.37%  "
$ws.Range("E26").Value = "  -4.43%  "
$ws.Range("D27").Value = "0.0₃0975"
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "20.50"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").Value = "162.52"
$ws.Range("E33").Value = "  +5.11%  "
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "4.76"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.541.77"
$ws.Range("E39").Value = "  -6.72%  "
$ws.Range("D40").Value = "22.82"
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.94"
$ws.Range("E41").Value = "  -3.84%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "37.70"
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("D44").Value = "0.0603"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0250"
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "5.10"
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "20.51"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "270.80"
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("E51").Value = "  +0.16%  "
